$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of OHLCV data to append (datetime serial, Open, High, Low, Close, Volume)
$data = @(
    @(45131.41666666666, 30091.5,   30102.82, 28874,    29189.1,  21487.26359591),
    @(45132.41666666666, 29188.72,  29378.38, 29065.48, 29236.54, 9918.333340929999),
    @(45133.41666666666, 29235.78,  29690.93, 29111.01, 29358.33, 13931.08744133),
    @(45134.41666666666, 29358,     29572.18, 29087.68, 29225.73, 10636.17094649),
    @(45135.41666666666, 29224.23,  29535.39, 29129.53, 29328.09, 9764.89835609),
    @(45136.41666666666, 29327.86,  29412.87, 29267.84, 29369.53, 4075.55845071),
    @(45137.41666666666, 29369.79,  29456.1,  29053.2,  29288.27, 6120.52992371),
    @(45138.41666666666, 29288.05,  29524.5,  29122.39, 29240.57, 11591.19804062),
    @(45139.41666666666, 29237.8,   29722.98, 28612,    29705.37, 24547.8278278),
    @(45140.41666666666, 29706.04,  30033.5,  28930.67, 29169.38, 24016.90813323)
)

$lastRow = 661
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]

    # Copy the formatting (style) of column A from the previous row so the
    # new datetime cells keep the same date number format / border / alignment
    $ws.Range("A" + $lastRow).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
